$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the original sheet, placing the copy right before it. This
# produces two identical sheets sharing the same data/shared-strings.
$ws1.Copy($ws1)

$copy = $wb.Worksheets.Item(1)
$orig = $wb.Worksheets.Item(2)

# Rename: the new copy becomes "Jailbroken Minecraft", the original
# (pre-existing) sheet becomes "Official Minecraft".
$copy.Name = "Jailbroken Minecraft"
$orig.Name = "Official Minecraft"

# The "Official Minecraft" (original) sheet is the one left active/selected.
$orig.Activate()
$orig.Select()

